$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.545.48'
$ws.Range('E2').Value = '  -0.75%  '
$ws.Range('D3').Value = '1.660.35'
$ws.Range('E3').Value = '  -3.56%  '
$ws.Range('E4').Value = '  +0.76%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.512'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.70%  '
$ws.Range('E7').Value = '  +0.77%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.29'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.23%  '
$ws.Range('E9').Value = '  -2.57%  '
$ws.Range('E10').Value = '  -1.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0875'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.32%  '
$ws.Range('D12').Value = '1.895.12'
$ws.Range('D13').Value = '1.660.53'
$ws.Range('E13').Value = '  -3.14%  '
$ws.Range('E14').Value = '  -2.89%  '
$ws.Range('E15').Value = '  -3.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.72'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '245.82'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.47%  '
$ws.Range('D18').Value = '27.542.51'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').Value = '0.0₃0730'
$ws.Range('E19').Value = '  -2.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.50%  '
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('E22').Value = '  -3.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.32%  '
$ws.Range('E24').Value = '  -4.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.99'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.17'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.21%  '
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('E28').Value = '  +0.78%  '
$ws.Range('E29').Value = '  -2.19%  '
$ws.Range('E30').Value = '  +4.96%  '
$ws.Range('E31').Value = '  -1.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.33'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.19%  '
$ws.Range('D33').Value = '1.438.89'
$ws.Range('E33').Value = '  -7.52%  '
$ws.Range('E34').Value = '  -5.43%  '
$ws.Range('E35').Value = '  -8.26%  '
$ws.Range('E36').Value = '  -0.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.928'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.579'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.86%  '
$ws.Range('E39').Value = '  -2.74%  '
$ws.Range('E40').Value = '  -2.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.05'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.25%  '
$ws.Range('E42').Value = '  +0.68%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.41'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.793'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('D45').Value = '1.803.02'
$ws.Range('E45').Value = '  -3.14%  '
$ws.Range('E46').Value = '  -3.62%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.53'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.48%  '
$ws.Range('E49').Value = '  -0.87%  '
$ws.Range('E50').Value = '  -4.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.79'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.73%  '
